$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) to be treated as text so that numeric-looking
# values (e.g. "587.21") are not silently converted to actual numbers by
# Excel's auto-detection. The Link/Coin/Volume columns never look like pure
# numbers so they do not need this treatment.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "60.410.73"
$ws.Range("E2").Value = "  +1.62%  "
$ws.Range("D3").Value = "2.605.00"
$ws.Range("E3").Value = "  +0.81%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "587.21"
$ws.Range("E5").Value = "  +6.39%  "
$ws.Range("D6").Value = "142.86"
$ws.Range("E6").Value = "  +1.37%  "
$ws.Range("D7").Value = "0.997"
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "0.600"
$ws.Range("E8").Value = "  +1.05%  "
$ws.Range("D9").Value = "2.613.24"
$ws.Range("E9").Value = "  +0.63%  "
$ws.Range("E10").Value = "  -3.03%  "
$ws.Range("D11").Value = "0.106"
$ws.Range("E11").Value = "  +1.60%  "
$ws.Range("D12").Value = "0.156"
$ws.Range("E12").Value = "  -3.15%  "
$ws.Range("D13").Value = "0.370"
$ws.Range("E13").Value = "  +4.69%  "
$ws.Range("D14").Value = "3.072.56"
$ws.Range("E14").Value = "  +0.97%  "
$ws.Range("D15").Value = "24.81"
$ws.Range("E15").Value = "  +6.96%  "
$ws.Range("D16").Value = "60.446.65"
$ws.Range("E16").Value = "  +1.68%  "
$ws.Range("E17").Value = "  +2.95%  "
$ws.Range("D18").Value = "2.612.73"
$ws.Range("E18").Value = "  +1.05%  "
$ws.Range("D19").Value = "11.37"
$ws.Range("E19").Value = "  +9.98%  "
$ws.Range("D20").Value = "4.66"
$ws.Range("E20").Value = "  +2.12%  "
$ws.Range("D21").Value = "347.40"
$ws.Range("E21").Value = "  +2.28%  "
$ws.Range("D22").Value = "6.92"
$ws.Range("E22").Value = "  +6.77%  "
$ws.Range("E23").Value = "  -0.18%  "
$ws.Range("D24").Value = "0.523"
$ws.Range("E24").Value = "  +9.60%  "
$ws.Range("D25").Value = "63.11"
$ws.Range("E25").Value = "  +0.13%  "
$ws.Range("D26").Value = "0.998"
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("D27").Value = "0.160"
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("D28").Value = "8.04"
$ws.Range("E28").Value = "  +7.94%  "
$ws.Range("D29").Value = "0.0₃0795"
$ws.Range("E29").Value = "  +3.26%  "
$ws.Range("E30").Value = "  +10.66%  "
$ws.Range("B31").Value = "USDe"
$ws.Range("C31").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D31").Value = "0.998"
$ws.Range("E31").Value = "  +0.10%  "
$ws.Range("B32").Value = "Aptos"
$ws.Range("C32").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D32").Value = "6.36"
$ws.Range("E32").Value = "  +3.34%  "
$ws.Range("D33").Value = "163.93"
$ws.Range("E33").Value = "  +3.75%  "
$ws.Range("D34").Value = "19.55"
$ws.Range("E34").Value = "  +2.01%  "
$ws.Range("D35").Value = "4.26"
$ws.Range("E35").Value = "  +3.24%  "
$ws.Range("E36").Value = "  +8.67%  "
$ws.Range("D37").Value = "1.23"
$ws.Range("E37").Value = "  +5.41%  "
$ws.Range("D38").Value = "1.63"
$ws.Range("E38").Value = "  +10.58%  "
$ws.Range("D39").Value = "37.94"
$ws.Range("E39").Value = "  +0.96%  "
$ws.Range("E40").Value = "  +6.58%  "
$ws.Range("D41").Value = "309.71"
$ws.Range("E41").Value = "  +7.76%  "
$ws.Range("E42").Value = "  -0.24%  "
$ws.Range("D43").Value = "134.97"
$ws.Range("E43").Value = "  -0.42%  "
$ws.Range("D44").Value = "0.0995"
$ws.Range("E44").Value = "  +2.29%  "
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("D46").Value = "19.78"
$ws.Range("E46").Value = "  +4.47%  "
$ws.Range("D47").Value = "0.605"
$ws.Range("E47").Value = "  +1.35%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "5.01"
$ws.Range("E48").Value = "  +5.70%  "
$ws.Range("B49").Value = "Hedera"
$ws.Range("C49").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D49").Value = "0.0550"
$ws.Range("E49").Value = "  +3.26%  "
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").Value = "20.19"
$ws.Range("E50").Value = "  +8.21%  "
$ws.Range("D51").Value = "0.0242"
$ws.Range("E51").Value = "  +3.40%  "

# Reset the style/format of the Price column back to the workbook default
# ("Normal") now that the text values are safely stored, so that no stray
# number-format styling is left behind on these cells.
$ws.Range("D2:D51").Style = "Normal"
